$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.4080523778287335
$ws.Range("C2").Value = -0.06805237782873351
$ws.Range("D2").Value = 0.2619476221712665
$ws.Range("E2").Value = 0.1798255795865583
$ws.Range("F2").Value = -0.1487002110073138
$ws.Range("G2").Value = 0.1519476221712665
$ws.Range("H2").Value = -0.0720523778287335
$ws.Range("B3").Value = -0.04304072301962425
$ws.Range("C3").Value = 0.2869592769803757
$ws.Range("D3").Value = 0.2048372343956676
$ws.Range("E3").Value = -0.1236885561982045
$ws.Range("F3").Value = 0.1769592769803758
$ws.Range("G3").Value = -0.04704072301962425
$ws.Range("B4").Value = 0.5227398812587823
$ws.Range("C4").Value = 0.4406178386740742
$ws.Range("D4").Value = 0.1120920480802021
$ws.Range("E4").Value = 0.4127398812587824
$ws.Range("F4").Value = 0.1887398812587824
$ws.Range("G4").Value = 0.3152499217222862
$ws.Range("H4").Value = -0.02551264069874778
$ws.Range("I4").Value = 0.2799154064837261
$ws.Range("J4").Value = 0.01851010070983469
$ws.Range("B5").Value = 1.413678501741462
$ws.Range("C5").Value = 1.08515271114759
$ws.Range("D5").Value = 1.38580054432617
$ws.Range("E5").Value = 1.16180054432617
$ws.Range("F5").Value = 1.288310584789674
$ws.Range("G5").Value = 0.9475480223686397
$ws.Range("H5").Value = 1.252976069551114
$ws.Range("I5").Value = 0.9915707637772222
$ws.Range("B6").Value = 0.1871668706385847
$ws.Range("C6").Value = 0.487814703817165
$ws.Range("D6").Value = 0.263814703817165
$ws.Range("E6").Value = 0.3903247442806688
$ws.Range("F6").Value = 0.04956218185963485
$ws.Range("G6").Value = 0.3549902290421088
$ws.Range("H6").Value = 0.09358492326821732
$ws.Range("B7").Value = 0.4237343126741154
$ws.Range("C7").Value = 0.1997343126741154
$ws.Range("D7").Value = 0.3262443531376191
$ws.Range("E7").Value = -0.01451820928341477
$ws.Range("F7").Value = 0.2909098378990592
$ws.Range("G7").Value = 0.02950453212516771
$ws.Range("B8").Value = 0.3642057909153072
$ws.Range("C8").Value = 0.4907158313788109
$ws.Range("D8").Value = 0.149953268957777
$ws.Range("E8").Value = 0.4553813161402509
$ws.Range("F8").Value = 0.1939760103663595
$ws.Range("G8").Value = 0.1687686335006937
$ws.Range("H8").Value = 0.3610269567322639
$ws.Range("I8").Value = 0.3122562552947308
$ws.Range("B9").Value = 0.304105594018323
$ws.Range("C9").Value = -0.03665696840271088
$ws.Range("D9").Value = 0.268771078779763
$ws.Range("E9").Value = 0.007365773005871601
$ws.Range("F9").Value = -0.01784160385979419
$ws.Range("G9").Value = 0.1744167193717761
$ws.Range("H9").Value = 0.1256460179342429
$ws.Range("B10").Value = -0.2336291894851882
$ws.Range("C10").Value = 0.07179885769728576
$ws.Range("D10").Value = -0.1896064480766057
$ws.Range("E10").Value = -0.2148138249422715
$ws.Range("F10").Value = -0.02255550171070122
$ws.Range("G10").Value = -0.07132620314823437
$ws.Range("B11").Value = 0.1620802390468526
$ws.Range("C11").Value = -0.09932506672703881
$ws.Range("D11").Value = -0.1245324435927046
$ws.Range("E11").Value = 0.06772587963886564
$ws.Range("F11").Value = 0.0189551782013325
$ws.Range("B12").Value = -0.237319901377073
$ws.Range("C12").Value = -0.2625272782427387
$ws.Range("D12").Value = -0.07026895501116853
$ws.Range("E12").Value = -0.1190396564487017
$ws.Range("B13").Value = -0.2357126085203348
$ws.Range("C13").Value = -0.04345428528876452
$ws.Range("D13").Value = -0.09222498672629767
$ws.Range("B14").Value = -0.1086967540443382
$ws.Range("C14").Value = -0.1574674554818714
$ws.Range("B15").Value = -0.08180776662314601
